# Fix minor typos in tech review presentation
#
# 1) Slide 4, "Content Placeholder 2": "Suppose webpage alread in 'table'
#    format..." -> "Suppose webpage already in 'table' format..."
# 2) Slide 8, "Content Placeholder 2": "Example) " -> "Example " (drop the
#    stray closing paren before the non-breaking space / link)

$p = $ppt.ActivePresentation

# --- Edit 1: slide 4 --------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$tr4 = $shp4.TextFrame.TextRange
$para4 = $tr4.Paragraphs(1, 1)
$fullRange4 = $tr4.Characters($para4.Start, $para4.Length)
$fullRange4.Text = "Suppose webpage already in 'table' format..."

# --- Edit 2: slide 8 ----------------------------------------------------
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(2)
$tr8 = $shp8.TextFrame.TextRange
$para8 = $tr8.Paragraphs(1, 1)
$nbsp = [char]0x00A0
$fixRange8 = $tr8.Characters($para8.Start, 9)
$fixRange8.Text = "Example" + $nbsp
